# "Generate Report for Handback" - update the localization-status report
# after a de-de / zh-cn handback: refresh the status text, fill in the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns for both language sheets, and widen a few columns that now hold
# longer file names.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- zh-cn handback columns ---
$mdFileName  = "94e35117-715b-43b5-9d4e-54bc1dad67bb.md"
$zhXlfName   = "94e35117-715b-43b5-9d4e-54bc1dad67bb.016a76dd21889117c16de60df1eb254461145ebb.zh-cn.xlf"
$deXlfName   = "94e35117-715b-43b5-9d4e-54bc1dad67bb.016a76dd21889117c16de60df1eb254461145ebb.de-de.xlf"

# Re-use the same target URL the existing "Source File Name" (A2) hyperlink
# already points at, for the new "Latest Target File" (I2) hyperlink.
$zhTargetUrl = $null
foreach ($h in $zhcn.Hyperlinks) { $zhTargetUrl = $h.Address }
$deTargetUrl = $null
foreach ($h in $dede.Hyperlinks) { $deTargetUrl = $h.Address }

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $zhTargetUrl, "", "", $mdFileName) | Out-Null
$zhcn.Range("J2").Value = $zhXlfName
$zhcn.Range("K2").Value = "2016-08-30 23:04:10"

# --- de-de handback columns ---
$dede.Hyperlinks.Add($dede.Range("I2"), $deTargetUrl, "", "", $mdFileName) | Out-Null
$dede.Range("J2").Value = $deXlfName
$dede.Range("K2").Value = "2016-08-30 23:04:18"

# --- widen columns that now hold the longer handback file names ---
$overview.Range("E1").ColumnWidth = 29.09
$overview.Range("F1").ColumnWidth = 29.09

$zhcn.Range("C1").ColumnWidth = 29.09
$zhcn.Range("I1").ColumnWidth = 39.09
$zhcn.Range("J1").ColumnWidth = 39.09

$dede.Range("C1").ColumnWidth = 29.09
$dede.Range("I1").ColumnWidth = 39.09
$dede.Range("J1").ColumnWidth = 39.09

Write-Host "Handback report updated."
